$d = $word.ActiveDocument

# Locate the paragraph that ends with "This is a Microsoft word document."
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pText = $d.Paragraphs($i).Range.Text
    if ($pText -like "This is a Microsoft word document.*") {
        $target = $d.Paragraphs($i)
        break
    }
}
if ($target -eq $null) {
    $target = $d.Paragraphs(1)
}

$pRange = $target.Range
$paraEnd = $pRange.End        # end of paragraph range (includes the paragraph mark)
$textEnd = $paraEnd - 1       # position right before the paragraph mark (end of visible text)

# Insert a new paragraph, immediately after the current one, that holds the
# three additional runs: " (", "Changed main", ")" -- each its own <w:r>,
# matching the structure produced by the target edit.
$insertionPoint = $d.Range($textEnd, $textEnd)
$xmlFragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
               '<w:r><w:t>Changed main</w:t></w:r>' +
               '<w:r><w:t>)</w:t></w:r>' +
               '</w:p>'
$insertionPoint.InsertXML($xmlFragment)

# The fragment above lands as a brand new paragraph right after the
# original one. Delete the original paragraph's mark so the new runs
# are pulled back up into the very same paragraph as the original run,
# rather than living in a paragraph of their own.
$mergeRange = $d.Range($textEnd, $textEnd + 1)
$mergeRange.Delete()
